# Update the public EPEX spot prices workbook
$wb = $excel.ActiveWorkbook

# --- "Prix Spot" sheet: insert a new date column (16-dec) before column EM ---
$ws = $wb.Worksheets.Item("Prix Spot")

# Insert a new column at position EM (column 143), shifting EM:FQ to FN:FR
$ws.Columns.Item(143).Insert()

# Header cell for the new column
$ws.Cells.Item(1, 143).Value = "16-dec"

# Fill the new column's data rows (2-25) with "-" like the surrounding empty cells
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 143).Value = "-"
}

# --- "Gaz" sheet: update the weekend-carried gas price values ---
$wsGaz = $wb.Worksheets.Item("Gaz")
$wsGaz.Range("B170").Value = 25.93
$wsGaz.Range("B171").Value = 25.93
